$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "20240430-20240506"
$ws.Range("B8").Value = "双周六道题，一直在看论文总结要点"
$ws.Range("D8").Value = "缺少灵感和创新点"

$ws.Rows.Item(8).RowHeight = 37

$ws.Range("D10").Select()
